$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the security/index labels in column A (rows 2-8)
$ws.Range("A2").Value = "PARSTEI LX Equity"
$ws.Range("A3").Value = "FLOT FP Equity"
$ws.Range("A4").Value = "SX5EEX GY Equity"
$ws.Range("A5").Value = "SPY US Equity"
$ws.Range("A6").Value = "LFGGBDR LX Equity"
$ws.Range("A7").Value = "EUN5 GY Equity"
$ws.Range("A8").Value = "EUNH GY Equity"

# Update the computed values in columns B (Opt Portfolio) and C (Opt Portfolio with View)
$ws.Range("C2").Value = [double]"0.9999999999999986"

$ws.Range("C3").Value = 0

$ws.Range("B4").Value = [double]"5.415120407071345e-16"
$ws.Range("C4").Value = [double]"8.292679237873184e-16"

$ws.Range("B5").Value = [double]"4.567492532905971e-16"
$ws.Range("C5").Value = [double]"6.684234565806446e-16"

$ws.Range("B6").Value = [double]"8.981681503671817e-16"
$ws.Range("C6").Value = [double]"1.077585560291168e-15"

$ws.Range("C7").Value = [double]"9.204482471635608e-16"

$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
